# Updated TEST to show usage of adate helpers
#
# This script:
#  1. Inserts a new "calculates" worksheet right after "survey" (new sheetId, becomes sheet3).
#  2. Populates its header row (calculation_name / calculation) with the matching font style.
#  3. Updates the "survey" sheet: changes the condition on row 6 and appends a new
#     if/else/end-if block (rows 11-16) that exercises the new adate.ageInYears /
#     adate.yearUnknown helpers.
#
# The order in which new string values are first assigned controls the order they
# receive in the shared string table, so that ordering is chosen deliberately to
# reproduce the target shared-strings layout.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")

# --- add the new "calculates" worksheet, right after "survey" ---
$calc = $wb.Worksheets.Add($null, $survey)
$calc.Name = "calculates"

# --- survey sheet edits (part 1, establishes shared-string order) ---
$survey.Range("F12").Value = "lblbla"
$survey.Range("F14").Value = "lblblabla"
$survey.Range("H12").Value = "Sand"
$survey.Range("H14").Value = "Falsk"

# --- calculates sheet header row ---
$calc.Range("A1").Value = "calculation_name"

$survey.Range("G12").Value = "Alder {{adate.ageInYears(data('ADA'))}}"
$survey.Range("C6").Value = "adate.ageInYears(data('ADA'))>2"
$survey.Range("G14").Value = "Ingen alder…"
$survey.Range("C11").Value = "!adate.yearUnknown(data('ADA'))"

$calc.Range("B1").Value = "calculation"

# remaining survey cells (reuse existing shared strings, order independent)
$survey.Range("B11").Value = "if"
$survey.Range("D12").Value = "note"
$survey.Range("B13").Value = "else"
$survey.Range("D14").Value = "note"
$survey.Range("B15").Value = "end if"
$survey.Range("B16").Value = "end screen"

# style the calculates header row with the new (grey) font color, matching cellXf 18 / fontId 8
$calc.Range("A1:B1").Font.Color = 4210752

# approximate the target column widths for the calculates sheet
$calc.Columns.Item(1).ColumnWidth = 15.83
$calc.Columns.Item(2).ColumnWidth = 25.5

$calc.PageSetup.Orientation = 1

# restore selections/active sheet to match the final authored state
$calc.Range("B1").Select() | Out-Null
$survey.Activate()
$survey.Range("C12").Select() | Out-Null
